# "added ch4 to GNFR C"
# GNFR sector "C" (row 3) previously had no CH4 [t] emissions recorded (F3 = 0).
# This edit adds the CH4 value for sector C and rolls the change into the
# CH4 [t] column total (row 16). The NOX (D3) and PM (G3) values for the
# same row shift by a single floating point ULP as a side effect of the
# upstream recalculation that produced the new figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sector "C" row (row 3)
$ws.Range("D3").Value = 610.0928791606466
$ws.Range("F3").Value = 70.69577479034739
$ws.Range("G3").Value = 4.094720970484224

# Recalculated column total for CH4 [t] (row 16)
$ws.Range("F16").Value = 1561.569662142907
